$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1943.031928535771
$ws.Range("B3").Value = 3330.891331388301
$ws.Range("B4").Value = 3955.303289966849
$ws.Range("B5").Value = 4315.689307757721
$ws.Range("B6").Value = 4377.759687877078
$ws.Range("B7").Value = 4678.418094215785
$ws.Range("B8").Value = 4729.39790750391
$ws.Range("B9").Value = 4720.455039049254
$ws.Range("B10").Value = 4684.478942018699
$ws.Range("B11").Value = 4393.626011160394
$ws.Range("B12").Value = 4306.883181863406
$ws.Range("B13").Value = 4178.109904109499
$ws.Range("B14").Value = 4080.951257608716
$ws.Range("B15").Value = 3929.829710870947
$ws.Range("B16").Value = 3538.958874814
$ws.Range("B17").Value = 3384.357685371022
$ws.Range("B18").Value = 3088.292263438086
$ws.Range("B19").Value = 2685.941369101534
$ws.Range("B20").Value = 2208.796493248246
$ws.Range("B21").Value = 2079.219694631881
$ws.Range("B22").Value = 1571.019483090939
$ws.Range("B23").Value = 1299.719302969215
$ws.Range("B24").Value = 891.6824229678857
$ws.Range("B25").Value = 488.4198119974228
$ws.Range("B26").Value = 238.2418642140564
$ws.Range("B27").Value = 76.0839726514977
$ws.Range("B28").Value = 76.0839726514977
$ws.Range("B29").Value = 76.0839726514977
$ws.Range("B30").Value = 76.0839726514977
$ws.Range("B31").Value = 76.0839726514977
$ws.Range("B32").Value = 76.0839726514977
$ws.Range("B33").Value = 76.0839726514977
$ws.Range("B34").Value = 76.0839726514977
$ws.Range("B35").Value = 76.0839726514977
$ws.Range("B36").Value = 76.0839726514977
$ws.Range("B37").Value = 76.0839726514977
$ws.Range("B38").Value = 76.0839726514977
$ws.Range("B39").Value = 76.0839726514977
$ws.Range("B40").Value = 76.0839726514977
$ws.Range("B41").Value = 76.0839726514977
$ws.Range("B42").Value = 76.0839726514977
$ws.Range("B43").Value = 76.0839726514977
$ws.Range("B44").Value = 76.0839726514977
$ws.Range("B45").Value = 76.0839726514977
$ws.Range("B46").Value = 76.0839726514977
$ws.Range("B47").Value = 76.0839726514977
$ws.Range("B48").Value = 76.0839726514977
$ws.Range("B49").Value = 76.0839726514977
$ws.Range("B50").Value = 76.0839726514977
$ws.Range("B51").Value = 76.0839726514977
$ws.Range("B52").Value = 76.0839726514977
$ws.Range("B53").Value = 76.0839726514977
$ws.Range("B54").Value = 76.0839726514977
$ws.Range("B55").Value = 76.0839726514977
$ws.Range("B56").Value = 76.0839726514977
$ws.Range("B57").Value = 76.0839726514977
$ws.Range("B58").Value = 76.0839726514977
$ws.Range("B59").Value = 76.0839726514977
$ws.Range("B60").Value = 76.0839726514977
$ws.Range("B61").Value = 76.0839726514977
$ws.Range("B62").Value = 76.0839726514977

